$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F2").Value = 73
$ws1.Range("F3").Value = 1073
$ws1.Range("F5").Value = 11129
$ws1.Range("F6").Value = 1289
$ws1.Range("F7").Value = 391
$ws1.Range("F8").Value = 630
$ws1.Range("F9").Value = 2072
$ws1.Range("F10").Value = 613
$ws1.Range("F11").Value = 777
$ws1.Range("F12").Value = 268
$ws1.Range("F13").Value = 365
$ws1.Range("F14").Value = 332
$ws1.Range("F15").Value = 309
$ws1.Range("F16").Value = 1035
$ws1.Range("F17").Value = 416
$ws1.Range("G17").Value = "'60"
$ws1.Range("F18").Value = 818
$ws1.Range("F19").Value = 247
$ws1.Range("F20").Value = 465
$ws1.Range("F21").Value = 705
$ws1.Range("F22").Value = 811
$ws1.Range("F23").Value = 30
$ws1.Range("F24").Value = 197
$ws1.Range("F25").Value = 496
$ws1.Range("F26").Value = 22
$ws1.Range("F27").Value = 9
$ws1.Range("F28").Value = 234

$ws2.Range("F4").Value = 655
$ws2.Range("F5").Value = 28
$ws2.Range("F6").Value = 63
$ws2.Range("F8").Value = 687
$ws2.Range("F9").Value = 12

$ws4.Range("F4").Value = 73
$ws4.Range("F5").Value = 1073
$ws4.Range("F6").Value = 655
$ws4.Range("F8").Value = 28
$ws4.Range("F9").Value = 11129
$ws4.Range("F10").Value = 1289
$ws4.Range("F11").Value = 63
$ws4.Range("F12").Value = 391
$ws4.Range("F13").Value = 630
$ws4.Range("F14").Value = 2072
$ws4.Range("F15").Value = 613
$ws4.Range("F16").Value = 777
$ws4.Range("F18").Value = 268
$ws4.Range("F19").Value = 365
$ws4.Range("F20").Value = 332
$ws4.Range("F21").Value = 309
$ws4.Range("F22").Value = 1035
$ws4.Range("F23").Value = 416
$ws4.Range("G23").Value = "'60"
$ws4.Range("F24").Value = 687
$ws4.Range("F25").Value = 818
$ws4.Range("F26").Value = 247
$ws4.Range("F27").Value = 465
$ws4.Range("F28").Value = 705
$ws4.Range("F29").Value = 811
$ws4.Range("F30").Value = 30
$ws4.Range("F31").Value = 12
$ws4.Range("F32").Value = 197
$ws4.Range("F33").Value = 496
$ws4.Range("F34").Value = 22
$ws4.Range("F35").Value = 9
$ws4.Range("F36").Value = 234
